$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": B11 4 -> 5, C11 -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 "Total": B12 88 -> 110, C12 -2 -> -2.4, E12 text "86/112" -> "107.6/140"
$ws.Range("B12").Value = 110
$ws.Range("C12").Value = -2.4
$ws.Range("E12").Value = "107.6/140"
